# Scheduled-runner price/profit refresh across the leve-crafting sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for
# the specific leves whose market data changed; some rows gain or lose a
# LeveProfitNQ/LeveProfitHQ (M/N) cell entirely when that recipe no longer
# has an HQ (or NQ) variant priced.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 73
$ws.Range("I9").Value = 64.5
$ws.Range("K9").Value = 64.5
$ws.Range("M9").Value = 104.5
$ws.Range("H17").Value = 3308
$ws.Range("J17").Value = 3308
$ws.Range("L17").Value = 9924
$ws.Range("N17").Value = -10260
$ws.Range("H19").Value = 949.5
$ws.Range("I19").Value = 900
$ws.Range("K19").Value = 900
$ws.Range("M19").Value = -725
$ws.Range("H32").Value = 1650
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 300
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 300
$ws.Range("M32").Value = -2674
$ws.Range("N32").Value = -952
$ws.Range("H62").Value = 6666.6665
$ws.Range("H65").Value = 6666.6665
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H112").Value = 4444
$ws.Range("J112").Value = 4444
$ws.Range("L112").Value = 13332
$ws.Range("N112").Value = -15548
$ws.Range("H137").Value = 3750
$ws.Range("J137").Value = 3500
$ws.Range("L137").Value = 10500
$ws.Range("N137").Value = -15600
$ws.Range("H138").Value = 5561.32
$ws.Range("J138").Value = 5751.0435
$ws.Range("L138").Value = 17253.1305
$ws.Range("N138").Value = -27533.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6849.8887
$ws.Range("I32").Value = 6076.353
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 6076.353
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -5789.353
$ws.Range("N32").Value = -20574
$ws.Range("H45").Value = 2690
$ws.Range("I45").Value = 2690
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2690
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2313
$ws.Range("N45").ClearContents()
$ws.Range("H102").Value = 2974
$ws.Range("I102").Value = 2949.5
$ws.Range("J102").Value = 2998.5
$ws.Range("K102").Value = 2949.5
$ws.Range("L102").Value = 2998.5
$ws.Range("M102").Value = -1327.5
$ws.Range("N102").Value = -6242.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H99").Value = 1200
$ws.Range("I99").Value = 1200
$ws.Range("K99").Value = 1200
$ws.Range("M99").Value = 298
$ws.Range("H105").Value = 3948
$ws.Range("I105").Value = 3948
$ws.Range("K105").Value = 3948
$ws.Range("M105").Value = -2201

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 70295
$ws.Range("J68").Value = 70295
$ws.Range("L68").Value = 70295
$ws.Range("N68").Value = -71793
$ws.Range("H71").Value = 70295
$ws.Range("J71").Value = 70295
$ws.Range("L71").Value = 210885
$ws.Range("N71").Value = -218373
$ws.Range("H132").Value = 1800.375
$ws.Range("I132").Value = 1376.1
$ws.Range("K132").Value = 4128.299999999999
$ws.Range("M132").Value = -1598.299999999999
$ws.Range("H141").Value = 99995
$ws.Range("J141").Value = 99995
$ws.Range("L141").Value = 99995
$ws.Range("N141").Value = -110355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 300
$ws.Range("M7").Value = -188
$ws.Range("H34").Value = 3749.0833
$ws.Range("J34").Value = 3999.182
$ws.Range("L34").Value = 11997.546
$ws.Range("N34").Value = -12165.546
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 799.5
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H122").Value = 5797
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 6996.25
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 62966.25
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -67866.25
$ws.Range("H129").Value = 1945.4
$ws.Range("I129").Value = 664.8333
$ws.Range("K129").Value = 1994.4999
$ws.Range("M129").Value = 3005.5001
$ws.Range("H131").Value = 2388.25
$ws.Range("I131").Value = 1715
$ws.Range("J131").Value = 2484.4285
$ws.Range("K131").Value = 5145
$ws.Range("L131").Value = 7453.2855
$ws.Range("M131").Value = -105
$ws.Range("N131").Value = -17533.2855
$ws.Range("H132").Value = 1440.5
$ws.Range("J132").Value = 999
$ws.Range("L132").Value = 8991
$ws.Range("N132").Value = -14051

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3090
$ws.Range("I132").Value = 2649.8
$ws.Range("K132").Value = 7949.400000000001
$ws.Range("M132").Value = -5419.400000000001
$ws.Range("H140").Value = 49982
$ws.Range("J140").Value = 49982
$ws.Range("L140").Value = 49982
$ws.Range("N140").Value = -60342

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 2000
$ws.Range("N27").Value = -2214
$ws.Range("H61").Value = 1749.5
$ws.Range("J61").Value = 2000
$ws.Range("L61").Value = 2000
$ws.Range("N61").Value = -2404
$ws.Range("H113").Value = 1749.5
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 849.5
$ws.Range("I107").Value = 849.5
$ws.Range("K107").Value = 2548.5
$ws.Range("M107").Value = -628.5
$ws.Range("H122").Value = 400
$ws.Range("I122").Value = 400
$ws.Range("K122").Value = 1200
$ws.Range("M122").Value = 1250
$ws.Range("H136").Value = 1242
$ws.Range("I136").Value = 1136.1666
$ws.Range("K136").Value = 3408.4998
$ws.Range("M136").Value = -858.4998000000001
